$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = "testing template"
$ws.Range("C68").Value = "testing template"
$ws.Range("D68").Value = 6
$ws.Range("E68").Value = "Custom"
$ws.Range("F68").Value = "testing_template"
$ws.Range("G68").Value = "custom"
$ws.Range("H68").Value = "template"
$ws.Range("I68").Value = "[]"
$ws.Range("K68").Value = "['node']"
$ws.Range("L68").Value = "template/"

# Row 69
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = "test 1"
$ws.Range("C69").Value = "lorem ipsum"
$ws.Range("D69").Value = 6
$ws.Range("E69").Value = "Custom"
$ws.Range("F69").Value = "test_1"
$ws.Range("G69").Value = "custom"
$ws.Range("H69").Value = "template"
$ws.Range("I69").Value = "[]"
$ws.Range("K69").Value = "['node']"
$ws.Range("L69").Value = "template/"
